$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 120, pushing existing rows 120:166 down to 121:167.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new record's data.
$ws.Cells.Item(120, 1).Value = 4
$ws.Cells.Item(120, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(120, 3).Value = "Los Lagos"
$ws.Cells.Item(120, 4).Value = 44524
$ws.Cells.Item(120, 5).Value = 10
$ws.Cells.Item(120, 6).Value = "Fruta"
$ws.Cells.Item(120, 7).Value = 100108
$ws.Cells.Item(120, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(120, 9).Value = 100108005
$ws.Cells.Item(120, 10).Value = "Piña"
$ws.Cells.Item(120, 11).Value = "Caramelo"
$ws.Cells.Item(120, 12).Value = "Tercera"
$ws.Cells.Item(120, 13).Value = 60
$ws.Cells.Item(120, 14).Value = 21000
$ws.Cells.Item(120, 15).Value = 22000
$ws.Cells.Item(120, 16).Value = 21500
$ws.Cells.Item(120, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(120, 18).Value = "Ecuador"
$ws.Cells.Item(120, 19).Value = 1344
$ws.Cells.Item(120, 20).Value = 16
